# Applies the edit described by the diff:
#  - Clear columns T:V for rows 12,15,16,17,18,19,31,34,35,36,37,38,43
#  - Change K46 from 5198000 to 50000
#  - Append 4 new data rows (47-50)
#  - Dimension grows to A1:V50 (handled automatically by the engine)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clear T:V for the rows whose extra (USD-converted) columns were removed.
#    NOTE: build the address with string concatenation, not "$r:" interpolation
#    -- "T$r:V$r" gets mis-parsed (the ":V$r" part is swallowed) by this
#    runtime's PowerShell-lite parser, which treats "$r:" like a scope prefix.
$rowsToClear = 12,15,16,17,18,19,31,34,35,36,37,38,43
foreach ($r in $rowsToClear) {
    $addr = "T" + $r + ":V" + $r
    $ws.Range($addr).ClearContents()
}

# 2) Fix the quantity typo on row 46
$ws.Range("K46").Value = 50000

# 3) Append the four new rows (47-50)
$newRows = @(
    @{ Row=47; A="import"; B="2019-08-01"; C="kakinada sea (inkak1)"; D="2018-05-25 00:00:00"; E=27131200; F="united states"; H="calcined petroleum coke"; I="usd"; J=485; K=100000; L="kgs"; M=33296.29; N=173074128.2; O="oxbow calcining international"; P="vedanta limited"; Q="panajigoa"; R="january"; S=2021; T=33296.29; U=173074128.2; V=485 },
    @{ Row=48; A="import"; B="2020-06-01"; C="kakinada sea (inkak1)"; D="2018-05-25 00:00:00"; E=27131200; F="united states"; H="calcined petroleum coke in bulk"; I="usd"; J=485; K=150000; L="kgs"; M=33296.29; N=173074128.2; O="oxbow calcining international"; P="vedanta limited"; Q="panajigoa"; R="january"; S=2021; T=33296.29; U=173074128.2; V=485 },
    @{ Row=49; A="import"; B="2021-07-01"; C="kakinada sea (inkak1)"; D="2018-05-25 00:00:00"; E=27131200; F="united states"; H="calcined petroleum coke in bulk"; I="usd"; J=485; K=180000; L="kgs"; M=33296.29; N=173074128.2; O="oxbow calcining international"; P="vedanta limited"; Q="panajigoa"; R="january"; S=2021; T=33296.29; U=173074128.2; V=485 },
    @{ Row=50; A="import"; B="2022-12-01"; C="kakinada sea (inkak1)"; D="2018-05-25 00:00:00"; E=27131200; F="united states"; H="calcined petroleum coke in bulk"; I="usd"; J=485; K=190000; L="kgs"; M=33296.29; N=173074128.2; O="oxbow calcining international"; P="vedanta limited"; Q="panajigoa"; R="january"; S=2021; T=33296.29; U=173074128.2; V=485 }
)

# Columns B looks like "YYYY-MM-DD" and the engine auto-coerces that to a
# real date serial; force it to stay literal text the same way a user
# typing a leading apostrophe would, then strip the resulting "Text" style
# so the cell comes back out with no explicit style (matching the diff).
$dateLikeCols = @("B")

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    foreach ($col in 'A','B','C','D','E','F','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V') {
        if ($rowData.ContainsKey($col)) {
            $addr = $col + $r
            $val = $rowData[$col]
            if ($dateLikeCols -contains $col) {
                $ws.Range($addr).Value = "'" + $val
                $ws.Range($addr).Style = "Normal"
            } else {
                $ws.Range($addr).Value = $val
            }
        }
    }
}
